$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data cells for rows 2-5 with new sensor readings
$ws.Range("A2").Value = 45061.50694444445
$ws.Range("B2").Value = 24.502
$ws.Range("C2").Value = 17.071
$ws.Range("D2").Value = 4.266
$ws.Range("E2").Value = 51.765
$ws.Range("F2").Value = 42.826
$ws.Range("G2").Value = 19.282
$ws.Range("H2").Value = 64.646
$ws.Range("I2").Value = 29.668
$ws.Range("J2").Value = 12.708
$ws.Range("K2").Value = 19.604
$ws.Range("L2").Value = 20.148
$ws.Range("M2").Value = 21.386
$ws.Range("N2").Value = 6.157
$ws.Range("O2").Value = 19.174
$ws.Range("P2").Value = 27.062
$ws.Range("Q2").Value = 15.864
$ws.Range("R2").Value = 3.991
$ws.Range("S2").Value = 2.607
$ws.Range("T2").Value = 284.514
$ws.Range("U2").Value = 53.424
$ws.Range("V2").Value = 17.698
$ws.Range("W2").Value = 35.645
$ws.Range("X2").Value = 18.491
$ws.Range("Y2").Value = 2.496
$ws.Range("Z2").Value = 31.956
$ws.Range("AA2").Value = 15.633
$ws.Range("AB2").Value = 13.958
$ws.Range("AC2").Value = 16.32
$ws.Range("AD2").Value = 21.117
$ws.Range("AE2").Value = 3.641
$ws.Range("AF2").Value = 57.291
$ws.Range("AG2").Value = 9.936999999999999
$ws.Range("AH2").Value = 22.127

$ws.Range("A3").Value = 45061.51388888889
$ws.Range("B3").Value = 24.021
$ws.Range("C3").Value = 17.38
$ws.Range("D3").Value = 2.069
$ws.Range("E3").Value = 51.658
$ws.Range("F3").Value = 42.771
$ws.Range("G3").Value = 18.904
$ws.Range("H3").Value = 73.45699999999999
$ws.Range("I3").Value = 29.086
$ws.Range("J3").Value = 12.763
$ws.Range("K3").Value = 19.309
$ws.Range("L3").Value = 20.534
$ws.Range("M3").Value = 21.75
$ws.Range("N3").Value = 6.039
$ws.Range("O3").Value = 18.798
$ws.Range("P3").Value = 26.681
$ws.Range("Q3").Value = 15.816
$ws.Range("R3").Value = 1.725
$ws.Range("S3").Value = 1.286
$ws.Range("T3").Value = 278.824
$ws.Range("U3").Value = 52.594
$ws.Range("V3").Value = 17.351
$ws.Range("W3").Value = 35.238
$ws.Range("X3").Value = 18.664
$ws.Range("Y3").Value = 2.447
$ws.Range("Z3").Value = 35.351
$ws.Range("AA3").Value = 15.326
$ws.Range("AB3").Value = 13.69
$ws.Range("AC3").Value = 16.055
$ws.Range("AD3").Value = 21.593
$ws.Range("AE3").Value = 1.294
$ws.Range("AF3").Value = 66.44499999999999
$ws.Range("AG3").Value = 9.769
$ws.Range("AH3").Value = 21.694

$ws.Range("A4").Value = 45061.52083333334
$ws.Range("B4").Value = 15.854
$ws.Range("C4").Value = 11.467
$ws.Range("D4").Value = 1.304
$ws.Range("E4").Value = 34.099
$ws.Range("F4").Value = 28.253
$ws.Range("G4").Value = 12.477
$ws.Range("H4").Value = 50.979
$ws.Range("I4").Value = 19.197
$ws.Range("J4").Value = 8.426
$ws.Range("K4").Value = 12.723
$ws.Range("L4").Value = 13.577
$ws.Range("M4").Value = 14.379
$ws.Range("N4").Value = 3.986
$ws.Range("O4").Value = 12.407
$ws.Range("P4").Value = 17.601
$ws.Range("Q4").Value = 10.513
$ws.Range("R4").Value = 1.095
$ws.Range("S4").Value = 0.792
$ws.Range("T4").Value = 181.534
$ws.Range("U4").Value = 34.765
$ws.Range("V4").Value = 11.452
$ws.Range("W4").Value = 23.249
$ws.Range("X4").Value = 12.385
$ws.Range("Y4").Value = 1.61
$ws.Range("Z4").Value = 24.063
$ws.Range("AA4").Value = 10.115
$ws.Range("AB4").Value = 9.058
$ws.Range("AC4").Value = 10.622
$ws.Range("AD4").Value = 14.311
$ws.Range("AE4").Value = 0.784
$ws.Range("AF4").Value = 46.099
$ws.Range("AG4").Value = 6.425
$ws.Range("AH4").Value = 14.318

$ws.Range("A5").Value = 45061.52777777778
$ws.Range("B5").Value = 15.37
$ws.Range("C5").Value = 11.22
$ws.Range("D5").Value = 1.07
$ws.Range("E5").Value = 33.15
$ws.Range("F5").Value = 27.46
$ws.Range("G5").Value = 12.1
$ws.Range("H5").Value = 47.89
$ws.Range("I5").Value = 18.62
$ws.Range("J5").Value = 8.210000000000001
$ws.Range("K5").Value = 12.36
$ws.Range("L5").Value = 13.23
$ws.Range("M5").Value = 14.01
$ws.Range("N5").Value = 3.86
$ws.Range("O5").Value = 12.03
$ws.Range("P5").Value = 17.09
$ws.Range("Q5").Value = 10.17
$ws.Range("R5").Value = 0.83
$ws.Range("S5").Value = 0.68
$ws.Range("T5").Value = 175.8
$ws.Range("U5").Value = 33.66
$ws.Range("V5").Value = 11.1
$ws.Range("W5").Value = 22.56
$ws.Range("X5").Value = 12.02
$ws.Range("Y5").Value = 1.56
$ws.Range("Z5").Value = 22.86
$ws.Range("AA5").Value = 9.81
$ws.Range("AB5").Value = 8.75
$ws.Range("AC5").Value = 10.27
$ws.Range("AD5").Value = 13.96
$ws.Range("AE5").Value = 0.5600000000000001
$ws.Range("AF5").Value = 43.26
$ws.Range("AG5").Value = 6.24
$ws.Range("AH5").Value = 13.88

# Remove the now-obsolete last data row (old row 6)
$ws.Rows.Item(6).Delete()

# Widen several data columns from 7 to 8 (and column T from 8 to 9) characters
$ws.Columns.Item(2).ColumnWidth = 7.17
$ws.Columns.Item(3).ColumnWidth = 7.17
$ws.Columns.Item(7).ColumnWidth = 7.17
$ws.Columns.Item(9).ColumnWidth = 7.17
$ws.Columns.Item(10).ColumnWidth = 7.17
$ws.Columns.Item(11).ColumnWidth = 7.17
$ws.Columns.Item(12).ColumnWidth = 7.17
$ws.Columns.Item(13).ColumnWidth = 7.17
$ws.Columns.Item(15).ColumnWidth = 7.17
$ws.Columns.Item(16).ColumnWidth = 7.17
$ws.Columns.Item(17).ColumnWidth = 7.17
$ws.Columns.Item(20).ColumnWidth = 8.17
$ws.Columns.Item(22).ColumnWidth = 7.17
$ws.Columns.Item(23).ColumnWidth = 7.17
$ws.Columns.Item(24).ColumnWidth = 7.17
$ws.Columns.Item(26).ColumnWidth = 7.17
$ws.Columns.Item(27).ColumnWidth = 7.17
$ws.Columns.Item(28).ColumnWidth = 7.17
$ws.Columns.Item(29).ColumnWidth = 7.17
$ws.Columns.Item(30).ColumnWidth = 7.17
$ws.Columns.Item(34).ColumnWidth = 7.17
